# Update aspirantes / candidates table:
#  - Row 2: replace Sandra Rodriguez Vega entry with Gabriel Perez Hernandez entry
#  - Row 3: replace Ana Karen Sandoval Ramos entry with Eduardo Ramos Guzman entry
#  - Row 4: replace Cesar Salinas Ramirez entry with Victor Campos Maurno entry
#  - Row 5 (new): add Hernandez Maria Martinez entry

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Cells.Item(2, 1).Value = "NI20200004"
$ws.Cells.Item(2, 3).Value = "PEHG921005HPLRRB04"
$ws.Cells.Item(2, 4).Value = "PEREZ"
$ws.Cells.Item(2, 5).Value = "HERNANDEZ"
$ws.Cells.Item(2, 6).Value = "GABRIEL"
$ws.Cells.Item(2, 7).Value = "GABRIEL PEREZ HERNANDEZ"
$ws.Cells.Item(2, 8).Value = "MASCULINO"
$ws.Cells.Item(2, 9).Value = "1993-03-08"
$ws.Cells.Item(2, 10).Value = "SÃ"

# --- Row 3 ---
$ws.Cells.Item(3, 1).Value = "NI20200005"
$ws.Cells.Item(3, 3).Value = "RAGE001211HPLMZDA0"
$ws.Cells.Item(3, 4).Value = "RAMOS"
$ws.Cells.Item(3, 5).Value = "GUZMAN"
$ws.Cells.Item(3, 6).Value = "EDUARDO"
$ws.Cells.Item(3, 7).Value = "EDUARDO RAMOS GUZMAN"
$ws.Cells.Item(3, 8).Value = "MASCULINO"
$ws.Cells.Item(3, 9).Value = "1993-03-09"

# --- Row 4 ---
$ws.Cells.Item(4, 1).Value = "NI20200007"
$ws.Cells.Item(4, 2).Value = "2020-09-01"
$ws.Cells.Item(4, 3).Value = "CAMV930125HPLMRC09"
$ws.Cells.Item(4, 4).Value = "CAMPOS "
$ws.Cells.Item(4, 5).Value = "MAURNO"
$ws.Cells.Item(4, 6).Value = "VICTOR"
$ws.Cells.Item(4, 7).Value = "VICTOR CAMPOS  MAURNO"
$ws.Cells.Item(4, 9).Value = "1993-03-11"

# --- Row 5 (new) ---
$ws.Cells.Item(5, 1).Value = "NI20210348"
$ws.Cells.Item(5, 2).Value = "2021-06-04"
$ws.Cells.Item(5, 3).Value = "MAHJ280603MSPRRV09"
$ws.Cells.Item(5, 4).Value = "MARÃA"
$ws.Cells.Item(5, 5).Value = "MARTÃNEZ"
$ws.Cells.Item(5, 6).Value = "HERNÃNDEZ"
$ws.Cells.Item(5, 7).Value = "HERNÃNDEZ MARÃA MARTÃNEZ"
$ws.Cells.Item(5, 8).Value = "FEMENINO"
$ws.Cells.Item(5, 9).Value = "1994-08-03"
$ws.Cells.Item(5, 10).Value = "NO"
$ws.Cells.Item(5, 11).Value = 28
$ws.Cells.Item(5, 12).Value = "a"
$ws.Cells.Item(5, 13).Value = "2021-06-04/2021-06-30"
